$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.682.69"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "3.684.66"
$ws.Range("E3").Value = "  +3.05%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "243.91"
$ws.Range("E5").Value = "  +1.42%  "
$ws.Range("E6").Value = "  +15.60%  "
$ws.Range("D7").Value = "668.22"
$ws.Range("E7").Value = "  +2.21%  "
$ws.Range("D8").Value = "0.429"
$ws.Range("E8").Value = "  +4.34%  "
$ws.Range("E9").Value = "  +4.26%  "
$ws.Range("E10").Value = "  +0.00%  "
$ws.Range("D11").Value = "3.679.67"
$ws.Range("E11").Value = "  +2.92%  "
$ws.Range("E12").Value = "  +4.72%  "
$ws.Range("E13").Value = "  +1.51%  "
$ws.Range("E14").Value = "  +3.31%  "
$ws.Range("D15").Value = "4.369.43"
$ws.Range("E15").Value = "  +3.08%  "
$ws.Range("D16").Value = "0.0000270"
$ws.Range("E16").Value = "  +3.80%  "
$ws.Range("D17").Value = "96.470.27"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").Value = "9.05"
$ws.Range("E18").Value = "  +16.47%  "
$ws.Range("D19").Value = "3.692.07"
$ws.Range("E19").Value = "  +3.40%  "
$ws.Range("D20").Value = "12.83"
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("D21").Value = "18.43"
$ws.Range("E21").Value = "  +3.45%  "
$ws.Range("E22").Value = "  +3.03%  "
$ws.Range("B23").Value = "SuiNetwork"
$ws.Range("C23").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D23").Value = "3.49"
$ws.Range("E23").Value = "  +2.54%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Value = "518.34"
$ws.Range("E24").Value = "  +3.19%  "
$ws.Range("E25").Value = "  +3.82%  "
$ws.Range("D26").Value = "6.97"
$ws.Range("E26").Value = "  +0.48%  "
$ws.Range("E27").Value = "  +6.76%  "
$ws.Range("D28").Value = "13.04"
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("E29").Value = "  +8.81%  "
$ws.Range("E30").Value = "  +1.41%  "
$ws.Range("D31").Value = "12.18"
$ws.Range("E31").Value = "  +6.85%  "
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("E33").Value = "  +1.87%  "
$ws.Range("D34").Value = "32.89"
$ws.Range("E34").Value = "  +4.86%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").Value = "1.81"
$ws.Range("E35").Value = "  +10.02%  "
$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D36").Value = "0.997"
$ws.Range("E36").Value = "  -0.65%  "
$ws.Range("E37").Value = "  +3.60%  "
$ws.Range("D38").Value = "619.93"
$ws.Range("E38").Value = "  -0.71%  "
$ws.Range("D39").Value = "8.75"
$ws.Range("E39").Value = "  -0.94%  "
$ws.Range("D40").Value = "42.59"
$ws.Range("E40").Value = "  +27.89%  "
$ws.Range("E41").Value = "  +6.32%  "
$ws.Range("E42").Value = "  +5.87%  "
$ws.Range("E43").Value = "  +7.15%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").Value = "6.18"
$ws.Range("E45").Value = "  +8.24%  "
$ws.Range("D46").Value = "0.0458"
$ws.Range("E47").Value = "  +25.38%  "
$ws.Range("D48").Value = "2.31"
$ws.Range("E48").Value = "  +1.26%  "
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("D50").Value = "8.61"
$ws.Range("E50").Value = "  +4.47%  "
$ws.Range("D51").Value = "54.59"
$ws.Range("E51").Value = "  +3.62%  "
